$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''62.884.00'
$ws.Range("E2").Value = '  -1.91%  '
$ws.Range("D3").Value = '''3.214.57'
$ws.Range("E3").Value = '  -2.57%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''527.17'
$ws.Range("E5").Value = '  +2.07%  '
$ws.Range("D6").Value = '''172.44'
$ws.Range("E6").Value = '  -4.69%  '
$ws.Range("D7").Value = '''0.592'
$ws.Range("E7").Value = '  +0.57%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '''3.211.47'
$ws.Range("E9").Value = '  -2.64%  '
$ws.Range("D10").Value = '''0.606'
$ws.Range("E10").Value = '  -1.19%  '
$ws.Range("D11").Value = '''53.22'
$ws.Range("E11").Value = '  -8.02%  '
$ws.Range("D12").Value = '''0.133'
$ws.Range("E12").Value = '  +2.09%  '
$ws.Range("D13").Value = '''0.0000252'
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("D14").Value = '''9.09'
$ws.Range("E14").Value = '  +0.90%  '
$ws.Range("D15").Value = '''3.741.79'
$ws.Range("E15").Value = '  -2.24%  '
$ws.Range("D16").Value = '''0.115'
$ws.Range("E16").Value = '  -4.54%  '
$ws.Range("D17").Value = '''3.220.66'
$ws.Range("E17").Value = '  -2.38%  '
$ws.Range("D18").Value = '''17.21'
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("D19").Value = '''62.910.91'
$ws.Range("E19").Value = '  -1.43%  '
$ws.Range("D20").Value = '''11.04'
$ws.Range("E20").Value = '  +2.52%  '
$ws.Range("D21").Value = '''0.966'
$ws.Range("E21").Value = '  +2.69%  '
$ws.Range("D22").Value = '''365.46'
$ws.Range("E22").Value = '  -1.11%  '
$ws.Range("D23").Value = '''3.77'
$ws.Range("E23").Value = '  +3.44%  '
$ws.Range("D24").Value = '''81.03'
$ws.Range("E24").Value = '  +1.67%  '
$ws.Range("D25").Value = '''11.02'
$ws.Range("E25").Value = '  +3.07%  '
$ws.Range("D26").Value = '''3.92'
$ws.Range("E26").Value = '  +5.11%  '
$ws.Range("E27").Value = '  +2.95%  '
$ws.Range("D28").Value = '''2.65'
$ws.Range("E28").Value = '  +0.79%  '
$ws.Range("D29").Value = '''11.29'
$ws.Range("E29").Value = '  +1.49%  '
$ws.Range("D30").Value = '''8.15'
$ws.Range("E30").Value = '  -1.11%  '
$ws.Range("D31").Value = '''28.46'
$ws.Range("E31").Value = '  +0.48%  '
$ws.Range("D32").Value = '''629.40'
$ws.Range("E32").Value = '  -2.40%  '
$ws.Range("D33").Value = '''6.46'
$ws.Range("E33").Value = '  -2.22%  '
$ws.Range("D34").Value = '''11.26'
$ws.Range("E34").Value = '  +2.13%  '
$ws.Range("E35").Value = '  +3.60%  '
$ws.Range("D36").Value = '''56.78'
$ws.Range("E36").Value = '  -3.93%  '
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").Value = '''36.74'
$ws.Range("E38").Value = '  +2.56%  '
$ws.Range("D39").Value = '''0.375'
$ws.Range("E39").Value = '  +1.03%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").Value = '''0.997'
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("B41").Value = 'PEPE'
$ws.Range("C41").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D41").Value = '''0.0₃0710'
$ws.Range("E41").Value = '  +14.89%  '
$ws.Range("D42").Value = '''0.123'
$ws.Range("E42").Value = '  +1.18%  '
$ws.Range("D43").Value = '''2.881.57'
$ws.Range("E43").Value = '  +3.58%  '
$ws.Range("D44").Value = '''2.51'
$ws.Range("E44").Value = '  +10.12%  '
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = '''2.68'
$ws.Range("E45").Value = '  +3.74%  '
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").Value = '''2.90'
$ws.Range("E46").Value = '  +10.44%  '
$ws.Range("D47").Value = '''0.0392'
$ws.Range("E47").Value = '  +2.49%  '
$ws.Range("D48").Value = '''2.59'
$ws.Range("E48").Value = '  -2.76%  '
$ws.Range("D49").Value = '''2.98'
$ws.Range("E49").Value = '  +8.58%  '
$ws.Range("D50").Value = '''0.124'
$ws.Range("E50").Value = '  +0.58%  '
$ws.Range("D51").Value = '''135.26'
$ws.Range("E51").Value = '  +0.93%  '
